$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows (2-33) with the new consumption / timestamp values
$ws.Cells.Item(2, 1).Value = 6566
$ws.Cells.Item(2, 2).Value = 46045.95833333334
$ws.Cells.Item(3, 1).Value = 6556
$ws.Cells.Item(3, 2).Value = 46045.96875
$ws.Cells.Item(4, 1).Value = 6434
$ws.Cells.Item(4, 2).Value = 46045.97916666666
$ws.Cells.Item(5, 1).Value = 6396
$ws.Cells.Item(5, 2).Value = 46045.98958333334
$ws.Cells.Item(6, 1).Value = 6320
$ws.Cells.Item(6, 2).Value = 46046
$ws.Cells.Item(7, 1).Value = 6219
$ws.Cells.Item(7, 2).Value = 46046.01041666666
$ws.Cells.Item(8, 1).Value = 6208
$ws.Cells.Item(8, 2).Value = 46046.02083333334
$ws.Cells.Item(9, 1).Value = 6156
$ws.Cells.Item(9, 2).Value = 46046.03125
$ws.Cells.Item(10, 1).Value = 6211
$ws.Cells.Item(10, 2).Value = 46046.04166666666
$ws.Cells.Item(11, 1).Value = 6160
$ws.Cells.Item(11, 2).Value = 46046.05208333334
$ws.Cells.Item(12, 1).Value = 6151
$ws.Cells.Item(12, 2).Value = 46046.0625
$ws.Cells.Item(13, 1).Value = 6123
$ws.Cells.Item(13, 2).Value = 46046.07291666666
$ws.Cells.Item(14, 1).Value = 6148
$ws.Cells.Item(14, 2).Value = 46046.08333333334
$ws.Cells.Item(15, 1).Value = 6111
$ws.Cells.Item(15, 2).Value = 46046.09375
$ws.Cells.Item(16, 1).Value = 6086
$ws.Cells.Item(16, 2).Value = 46046.10416666666
$ws.Cells.Item(17, 1).Value = 6106
$ws.Cells.Item(17, 2).Value = 46046.11458333334
$ws.Cells.Item(18, 1).Value = 6133
$ws.Cells.Item(18, 2).Value = 46046.125
$ws.Cells.Item(19, 1).Value = 6108
$ws.Cells.Item(19, 2).Value = 46046.14583333334
$ws.Cells.Item(20, 1).Value = 6139
$ws.Cells.Item(20, 2).Value = 46046.15625
$ws.Cells.Item(21, 1).Value = 6164
$ws.Cells.Item(21, 2).Value = 46046.16666666666
$ws.Cells.Item(22, 1).Value = 6166
$ws.Cells.Item(22, 2).Value = 46046.17708333334
$ws.Cells.Item(23, 1).Value = 6200
$ws.Cells.Item(23, 2).Value = 46046.1875
$ws.Cells.Item(24, 1).Value = 6242
$ws.Cells.Item(24, 2).Value = 46046.19791666666
$ws.Cells.Item(25, 1).Value = 6259
$ws.Cells.Item(25, 2).Value = 46046.20833333334
$ws.Cells.Item(26, 1).Value = 6320
$ws.Cells.Item(26, 2).Value = 46046.21875
$ws.Cells.Item(27, 1).Value = 6329
$ws.Cells.Item(27, 2).Value = 46046.22916666666
$ws.Cells.Item(28, 1).Value = 6411
$ws.Cells.Item(28, 2).Value = 46046.23958333334
$ws.Cells.Item(29, 1).Value = 6536
$ws.Cells.Item(29, 2).Value = 46046.25
$ws.Cells.Item(30, 1).Value = 6619
$ws.Cells.Item(30, 2).Value = 46046.26041666666
$ws.Cells.Item(31, 1).Value = 6697
$ws.Cells.Item(31, 2).Value = 46046.27083333334
$ws.Cells.Item(32, 1).Value = 6805
$ws.Cells.Item(32, 2).Value = 46046.28125
$ws.Cells.Item(33, 1).Value = 6907
$ws.Cells.Item(33, 2).Value = 46046.29166666666

# The last two data rows (34 and 35) are no longer present in the refreshed
# dataset, so remove them entirely (shrinks dimension from A1:B35 to A1:B33)
$ws.Range("A34:B35").Delete()
